$d = $word.ActiveDocument

# 1) Delete the first table (Имя/Sara, Фамилия/Blade222222, Пол/Ж, Возраст/53, Доход/167)
$d.Tables.Item(1).Delete()

# 2) Delete the standalone "Jhonson" paragraph (keep the "Sara" paragraph before it)
$rng = $d.Content
$found = $rng.Find.Execute("Jhonson")
$para = $rng.Paragraphs.Item(1)
$delRange = $d.Range($para.Range.Start, $para.Range.End + 1)
$delRange.Delete()

# 3) Update the remaining table's values: Jhon/Wick/М/33/345,7 -> Sara/Blade222222/Ж/53/167
$t = $d.Tables.Item(1)
$t.Cell(1, 2).Range.Text = "Sara"
$t.Cell(2, 2).Range.Text = "Blade222222"
$t.Cell(3, 2).Range.Text = "Ж"
$t.Cell(4, 2).Range.Text = "53"
$t.Cell(5, 2).Range.Text = "167"

# 4) Rename the trailing "Jhon" paragraph (after the table) to "Jhonson"
$searchStart = $t.Range.End
$rngTail = $d.Range($searchStart, $d.Content.End)
$foundTail = $rngTail.Find.Execute("Jhon")
$rngTail.Text = "Jhonson"

# 5) Insert a new "Jack" paragraph right after the "Jhonson" paragraph
$insertPos = $d.Content.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$jackParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Jack</w:t></w:r></w:p>
'@
$insertRange.InsertXML($jackParaXml)

# 6) Insert the new table (Имя/Jack, Фамилия/Daniels1111111, Пол/М, Возраст/21, Доход/12,5)
#    right after the "Jack" paragraph, before the trailing empty paragraph.
$tblInsertPos = $d.Content.End - 1
$tblInsertRange = $d.Range($tblInsertPos, $tblInsertPos)
$newTableXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="a3"/><w:tblW w:w="0" w:type="auto"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4677"/><w:gridCol w:w="4678"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="4677" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Имя</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4678" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Jack</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4677" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Фамилия</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4678" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Daniels1111111</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4677" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Пол</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4678" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>М</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4677" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Возраст</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4678" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>21</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="4677" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Доход</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4678" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>12,5</w:t></w:r></w:p></w:tc></w:tr></w:tbl>
'@
$tblInsertRange.InsertXML($newTableXml)

Write-Host "Edit complete"
